$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$values = @{
    2  = "-2,-9"
    3  = "8,-6"
    4  = "-7,-10"
    5  = "-8,3"
    6  = "-8,5"
    7  = "-7,-7"
    8  = "-1,2"
    9  = "-7,-9"
    10 = "-8,0"
    11 = "3,-3"
    12 = "2,3"
    13 = "6,9"
    14 = "5,4"
    15 = "-4,4"
    16 = "-5,3"
    17 = "4,6"
    18 = "-7,3"
    19 = "-7,0"
    20 = "4,-1"
    21 = "-6,8"
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
